$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width change (stored xml width 12.7109375 -> 14.7109375, i.e. +2 characters).
# The COM ColumnWidth->stored-width conversion in this host only lands on
# sixth-of-a-character steps, so 13.83 is the input that rounds to the
# closest achievable stored width (14.666...) to the target 14.7109375.
$ws.Range("E1").ColumnWidth = 13.83

# Row 2
$ws.Range("B2").Value = 0.027006397956728607
$ws.Range("C2").Value = 0.026893981217712914
$ws.Range("D2").Value = 0.026903019585067682
$ws.Range("E2").Value = 0.004364121870257579

# Row 3
$ws.Range("B3").Value = 30.019776248350702
$ws.Range("C3").Value = 30.018817954391618
$ws.Range("D3").Value = 30.011538791421234
$ws.Range("E3").Value = 99.956010938252788

# Row 4
$ws.Range("B4").Value = 0.1074643643580597
$ws.Range("C4").Value = 0.10529127196098649
$ws.Range("D4").Value = 0.10551680732681155
$ws.Range("E4").Value = 0.10940963220711369

# Row 5
$ws.Range("B5").Value = 1.8299370463193918
$ws.Range("C5").Value = 1.84734396868327
$ws.Range("D5").Value = 1.8453924521164513
$ws.Range("E5").Value = 0.9965561804803903

# Row 6
$ws.Range("E6").Value = 0.013734688769642182

# Row 9
$ws.Range("B9").Value = 249.95797553755358
$ws.Range("C9").Value = 249.98667448699553
$ws.Range("D9").Value = 249.98836639623838
$ws.Range("E9").Value = 150.27676897481228

# Row 10
$ws.Range("E10").Value = 17.676729830555718

# Row 12
$ws.Range("B12").Value = 25.941207743364313
$ws.Range("C12").Value = 26.24061612171343
$ws.Range("D12").Value = 26.253627802783942
$ws.Range("E12").Value = 7.619290605319506

# Row 15
$ws.Range("B15").Value = 0.010734176645309927
$ws.Range("C15").Value = 0.010606216411124881
$ws.Range("D15").Value = 0.01062345372257986
$ws.Range("E15").Value = 0.0041618191671662684

# Row 16
$ws.Range("B16").Value = 0.66435766594882084
$ws.Range("C16").Value = 0.66138126033990485
$ws.Range("D16").Value = 0.66180835950768879
$ws.Range("E16").Value = 0.46492856340000921
